# fall 13 week 3 inputs
# Append 23 new matchup rows (1143-1165) to the "Nine" sheet and move the
# viewport/selection down to follow them, mirroring what Excel does after a
# manual paste at the bottom of a growing data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(6, 18, 5, 2),
    @(4, 7, 3, 13),
    @(4, 8, 3, 12),
    @(2, 13, 3, 7),
    @(5, 4, 4, 16),
    @(4, 7, 2, 13),
    @(4, 4, 5, 16),
    @(4, 13, 6, 7),
    @(3, 5, 2, 15),
    @(5, 4, 3, 16),
    @(2, 14, 3, 6),
    @(4, 5, 2, 15),
    @(5, 13, 4, 7),
    @(5, 7, 9, 13),
    @(5, 2, 4, 18),
    @(4, 5, 3, 15),
    @(3, 16, 4, 4),
    @(7, 5, 4, 15),
    @(5, 12, 7, 8),
    @(3, 15, 2, 5),
    @(2, 15, 6, 5),
    @(5, 20, 2, 0),
    @(3, 14, 4, 6)
)

$firstNewRow = 1143
$lastNewRow = $firstNewRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Scroll the window so row 1148 is at the top (best effort - not every host
# persists the unfrozen scroll position, but harmless either way) and select
# the cell right after the freshly-entered block, as Excel leaves the active
# cell there after a block of manual entry.
$nextCell = $ws.Cells.Item($lastNewRow + 1, 1)
try {
    $excel.ActiveWindow.ScrollRow = 1148
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$nextCell.Select()
